$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "transt0"
$ws.Range("B2").Value = "Rút tiền"
$ws.Range("C2").Value = "10/11/2024 12:00:00 SA"
$ws.Range("E2").Value = "aksldfjalsdf"

# "99999" looks numeric, so force it to be stored as text (matching the
# source data which keeps it as a shared string) by entering it with a
# leading quote prefix, then restoring the cell's original formatting.
$ws.Range("D2").Value = "'99999"
$ws.Range("A1").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
